$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GOCharacterStatTable")

# Revert "Test: Skill Set in CharacterDataTable" — remove the SkillQ/SkillW/SkillE/SkillR
# columns (B:E) that were inserted after CHARACTER_NAME, restoring the original layout.
$ws.Range("B:E").EntireColumn.Delete()

# Widen the now-merged name/base-skill columns back to their pre-test sizes.
$ws.Range("A:B").ColumnWidth = 29.571428571428573
$ws.Range("C:C").ColumnWidth = 22.571428571428573

$ws.Range("F7").Select() | Out-Null
